# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.454.40"
$ws.Range("E2").Value = "  -2.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.361.01"
$ws.Range("E3").Value = "  -4.26%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.14"
$ws.Range("E5").Value = "  -4.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.99"
$ws.Range("E6").Value = "  -7.71%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.358.87"
$ws.Range("E8").Value = "  -4.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("E9").Value = "  -3.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  -4.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").Value = "  -4.79%  "

$ws.Range("E12").Value = "  -4.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.945.96"
$ws.Range("E13").Value = "  -3.92%  "

$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.370.41"
$ws.Range("E15").Value = "  -3.94%  "

$ws.Range("E16").Value = "  -6.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.536.02"
$ws.Range("E17").Value = "  -2.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.28"
$ws.Range("E18").Value = "  -5.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.28"
$ws.Range("E19").Value = "  -7.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.58"
$ws.Range("E20").Value = "  -3.14%  "

$ws.Range("E21").Value = "  -4.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "368.29"
$ws.Range("E22").Value = "  -6.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.501.78"
$ws.Range("E24").Value = "  -4.10%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.95"
$ws.Range("E26").Value = "  -4.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000104"
$ws.Range("E27").Value = "  -11.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  -6.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("E30").Value = "  -7.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.69"
$ws.Range("E31").Value = "  -7.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.148"
$ws.Range("E33").Value = "  -5.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.396.11"
$ws.Range("E34").Value = "  -4.05%  "

$ws.Range("E35").Value = "  -7.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.54"
$ws.Range("E36").Value = "  -3.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.08"
$ws.Range("E37").Value = "  -4.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "164.51"
$ws.Range("E38").Value = "  -1.63%  "

$ws.Range("E39").Value = "  -5.39%  "

$ws.Range("E40").Value = "  -6.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0742"
$ws.Range("E41").Value = "  -5.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.764"
$ws.Range("E43").Value = "  -5.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.22"
$ws.Range("E44").Value = "  -2.07%  "

$ws.Range("E45").Value = "  -5.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.48"
$ws.Range("E46").Value = "  -9.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.52"
$ws.Range("E47").Value = "  -8.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.05"
$ws.Range("E48").Value = "  -10.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.58"
$ws.Range("E49").Value = "  -3.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.227.65"
$ws.Range("E50").Value = "  -6.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.839"
$ws.Range("E51").Value = "  -6.63%  "
